$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 19:29"

# --- Refresh per-country stats (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos (row 4)
$ws.Cells.Item(4,2).Value = 6270129
$ws.Cells.Item(4,3).Value = 12558
$ws.Cells.Item(4,4).Value = 3506751
$ws.Cells.Item(4,5).Value = 2574127
$ws.Cells.Item(4,7).Value = 351
$ws.Cells.Item(4,8).Value = 189251

# Brasil (row 5)
$ws.Cells.Item(5,2).Value = 3961502
$ws.Cells.Item(5,3).Value = 8712
$ws.Cells.Item(5,5).Value = 679465
$ws.Cells.Item(5,7).Value = 260
$ws.Cells.Item(5,8).Value = 122941

# India (row 6)
$ws.Cells.Item(6,2).Value = 3826387
$ws.Cells.Item(6,3).Value = 60279
$ws.Cells.Item(6,4).Value = 2950122
$ws.Cells.Item(6,5).Value = 809082
$ws.Cells.Item(6,7).Value = 723
$ws.Cells.Item(6,8).Value = 67183

# Espana (row 12)
$ws.Cells.Item(12,2).Value = 479554
$ws.Cells.Item(12,3).Value = 8581
$ws.Cells.Item(12,7).Value = 42
$ws.Cells.Item(12,8).Value = 29194

# Chile (row 14)
$ws.Cells.Item(14,2).Value = 414739
$ws.Cells.Item(14,3).Value = 1594
$ws.Cells.Item(14,4).Value = 387683
$ws.Cells.Item(14,5).Value = 15712
$ws.Cells.Item(14,7).Value = 23
$ws.Cells.Item(14,8).Value = 11344

# Francia (row 20)
$ws.Cells.Item(20,2).Value = 293024
$ws.Cells.Item(20,3).Value = 7017
$ws.Cells.Item(20,5).Value = 175651

# Alemania (row 23)
$ws.Cells.Item(23,2).Value = 246808
$ws.Cells.Item(23,3).Value = 807
$ws.Cells.Item(23,5).Value = 15623

# Israel (row 29)
$ws.Cells.Item(29,2).Value = 121023
$ws.Cells.Item(29,3).Value = 2485
$ws.Cells.Item(29,4).Value = 97218
$ws.Cells.Item(29,5).Value = 22836
$ws.Cells.Item(29,7).Value = 12
$ws.Cells.Item(29,8).Value = 969

# Catar (row 30)
$ws.Cells.Item(30,2).Value = 119206
$ws.Cells.Item(30,3).Value = 212
$ws.Cells.Item(30,4).Value = 116111
$ws.Cells.Item(30,5).Value = 2896
$ws.Cells.Item(30,7).Value = 1
$ws.Cells.Item(30,8).Value = 199

# Ecuador (row 32)
$ws.Cells.Item(32,2).Value = 115457
$ws.Cells.Item(32,3).Value = 1148
$ws.Cells.Item(32,4).Value = 101867
$ws.Cells.Item(32,5).Value = 6971
$ws.Cells.Item(32,7).Value = 48
$ws.Cells.Item(32,8).Value = 6619

# Republica Dominicana (row 35)
$ws.Cells.Item(35,2).Value = 95627
$ws.Cells.Item(35,3).Value = 648
$ws.Cells.Item(35,4).Value = 69519
$ws.Cells.Item(35,5).Value = 24343
$ws.Cells.Item(35,7).Value = 27
$ws.Cells.Item(35,8).Value = 1765

# Chequia (row 74)
$ws.Cells.Item(74,2).Value = 25579
$ws.Cells.Item(74,3).Value = 462
$ws.Cells.Item(74,4).Value = 18307
$ws.Cells.Item(74,5).Value = 6847

# --- Libano moves up the ranking (new figures), Paraguay and Costa de
#     Marfil each shift down one row as a result ---

# Row 79 becomes Libano with fresh figures
$ws.Cells.Item(79,1).Value = "Libano"
$ws.Cells.Item(79,2).Value = 18375
$ws.Cells.Item(79,3).Value = 598
$ws.Cells.Item(79,4).Value = 5195
$ws.Cells.Item(79,5).Value = 13003
$ws.Cells.Item(79,7).Value = 6
$ws.Cells.Item(79,8).Value = 177

# Row 80 becomes Paraguay (its previous, still-current figures)
$ws.Cells.Item(80,1).Value = "Paraguay"
$ws.Cells.Item(80,2).Value = 18338
$ws.Cells.Item(80,3).Value = 0
$ws.Cells.Item(80,4).Value = 9721
$ws.Cells.Item(80,5).Value = 8269
$ws.Cells.Item(80,8).Value = 348

# Row 81 becomes Costa de Marfil (its previous, still-current figures)
$ws.Cells.Item(81,1).Value = "Costa de Marfil"
$ws.Cells.Item(81,2).Value = 18103
$ws.Cells.Item(81,3).Value = 0
$ws.Cells.Item(81,4).Value = 16814
$ws.Cells.Item(81,5).Value = 1172
$ws.Cells.Item(81,8).Value = 117

# Maldivas (row 102)
$ws.Cells.Item(102,2).Value = 8140
$ws.Cells.Item(102,3).Value = 137
$ws.Cells.Item(102,4).Value = 5338
$ws.Cells.Item(102,5).Value = 2773

# Mali (row 133)
$ws.Cells.Item(133,2).Value = 2802
$ws.Cells.Item(133,3).Value = 25
$ws.Cells.Item(133,4).Value = 2185
$ws.Cells.Item(133,5).Value = 491

# Sudan del Sur (row 136)
$ws.Cells.Item(136,2).Value = 2532
$ws.Cells.Item(136,3).Value = 5
$ws.Cells.Item(136,5).Value = 1195
